$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Servos" (sheet1): redesign the lift to use two slats (arms) instead
# of a single lever + separate "beams for mount cover"/chain bookkeeping.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Servos")

# --- Surgical edits: touch only what actually changed, so that untouched
# formulas (B7, B8) keep their original textual form. ---

# Drop the "beams for mount cover" column (G) that is no longer used.
$ws1.Range("G1").ClearContents()
$ws1.Range("G2").ClearContents()
$ws1.Range("G3").ClearContents()
$ws1.Range("G4").ClearContents()

# The lever length for the bucket/cover servos is renamed to "arm of force".
# (A4's new string is registered before A3's so the shared-string table
# ends up in the same order as the target workbook.)
$ws1.Range("A4").Value = "arm of force for cover servo,cm"
$ws1.Range("A3").Value = "arm of force for bucket servo,cm"

# The cover arm is now 2 cm (was 14, using the old single-lever geometry).
$ws1.Range("F4").Value = 2

# The "lever for bucket-to-side turn" (E6) moves to J2 (two-slat offset).
$ws1.Range("E6").ClearContents()
$ws1.Range("J2").Value = 4.8

# New row 5: arm of force for the servo that turns the bucket to the side,
# now computed off row 3 (bucket arm) plus the two-slat offset J2, for every
# column (including G, which used to be a "shared" formula off row 12).
$ws1.Range("A5").Value = "arm of force for servo that turn bucket to side"
$ws1.Range("B5").Formula = '=B3+$J$2'
$ws1.Range("C5").Formula = '=C3+$J$2'
$ws1.Range("D5").Formula = '=D3+$J$2'
$ws1.Range("E5").Formula = '=E3+$J$2'
$ws1.Range("F5").Formula = '=F3+$J$2'
$ws1.Range("G5").Formula = '=G3+$J$2'
# D5 references D3 (which carries a "0.00" number format); avoid leaking
# that format onto D5, which stays General like the rest of the row.
$ws1.Range("D5").NumberFormat = "General"

# Rename "bucket servo" (vertical move) to "bucket up-down servo".
$ws1.Range("A7").Value = "bucket up-down servo"

# J8 takes over the old E10 "weight" input for the side-turn servo torque.
$ws1.Range("J8").Value = 70

# New row 9: "bucket to side servo" torque, replacing the old A10/E10 pair;
# references row 5 (arm of force) and J8 (weight) instead of row 12/E6/E10.
$ws1.Range("A9").Value = "bucket to side servo"
$ws1.Range("B9").Formula = '=(B2*B5+C2*C5+D2*D5+E2*E5+F2*F5+G2*G5 + J2*J8) * H2 / 1000'
# B9 indirectly references D5/D3 (the "0.00" column); keep it General too.
$ws1.Range("B9").NumberFormat = "General"

# Remove the old row 10 (bucket-to-side torque) and row 12 (shared lever
# formulas) - their content has been superseded by rows 5 and 9 above.
$ws1.Range("A10").ClearContents()
$ws1.Range("E10").ClearContents()
$ws1.Range("A12:G12").ClearContents()

$ws1.Activate()
$ws1.Range("E11").Select()

# ---------------------------------------------------------------------------
# Sheet "Motor" (sheet2): drop the "chain" column (old column F) entirely -
# deleting the column shifts G/H/I left into F/G/H and K into J, carrying
# all the labels/values along automatically.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Motor")
$ws2.Columns("F:F").Delete()

# The torque formula drops the old "chain" term (I2*I3) and now multiplies
# by J2 (shifted from K2).
$ws2.Range("B6").Formula = '=(B2*B3+C2*C3+D2*D3+E2*E3+F2*F3+G2*G3+H2*H3) * J2 / 1000'

$ws2.Range("C8").Select()
